$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Formula = '="' + '66.198.82' + '"'
$ws.Range("D2").Copy()
$ws.Range("D2").PasteSpecial(-4163)

$ws.Range("D3").Formula = '="' + '3.320.80' + '"'
$ws.Range("D3").Copy()
$ws.Range("D3").PasteSpecial(-4163)

$ws.Range("E3").Formula = '="' + '  +0.67%  ' + '"'
$ws.Range("E3").Copy()
$ws.Range("E3").PasteSpecial(-4163)

$ws.Range("E4").Formula = '="' + '  -0.11%  ' + '"'
$ws.Range("E4").Copy()
$ws.Range("E4").PasteSpecial(-4163)

$ws.Range("D5").Formula = '="' + '564.31' + '"'
$ws.Range("D5").Copy()
$ws.Range("D5").PasteSpecial(-4163)

$ws.Range("E5").Formula = '="' + '  +1.31%  ' + '"'
$ws.Range("E5").Copy()
$ws.Range("E5").PasteSpecial(-4163)

$ws.Range("D6").Formula = '="' + '186.02' + '"'
$ws.Range("D6").Copy()
$ws.Range("D6").PasteSpecial(-4163)

$ws.Range("E6").Formula = '="' + '  +1.56%  ' + '"'
$ws.Range("E6").Copy()
$ws.Range("E6").PasteSpecial(-4163)

$ws.Range("E7").Formula = '="' + '  +0.02%  ' + '"'
$ws.Range("E7").Copy()
$ws.Range("E7").PasteSpecial(-4163)

$ws.Range("D8").Formula = '="' + '3.314.63' + '"'
$ws.Range("D8").Copy()
$ws.Range("D8").PasteSpecial(-4163)

$ws.Range("E8").Formula = '="' + '  +0.57%  ' + '"'
$ws.Range("E8").Copy()
$ws.Range("E8").PasteSpecial(-4163)

$ws.Range("E9").Formula = '="' + '  -2.00%  ' + '"'
$ws.Range("E9").Copy()
$ws.Range("E9").PasteSpecial(-4163)

$ws.Range("E10").Formula = '="' + '  -4.81%  ' + '"'
$ws.Range("E10").Copy()
$ws.Range("E10").PasteSpecial(-4163)

$ws.Range("D11").Formula = '="' + '0.574' + '"'
$ws.Range("D11").Copy()
$ws.Range("D11").PasteSpecial(-4163)

$ws.Range("E11").Formula = '="' + '  -1.79%  ' + '"'
$ws.Range("E11").Copy()
$ws.Range("E11").PasteSpecial(-4163)

$ws.Range("D12").Formula = '="' + '46.12' + '"'
$ws.Range("D12").Copy()
$ws.Range("D12").PasteSpecial(-4163)

$ws.Range("E12").Formula = '="' + '  -2.65%  ' + '"'
$ws.Range("E12").Copy()
$ws.Range("E12").PasteSpecial(-4163)

$ws.Range("D13").Formula = '="' + '0.0000265' + '"'
$ws.Range("D13").Copy()
$ws.Range("D13").PasteSpecial(-4163)

$ws.Range("E13").Formula = '="' + '  -0.36%  ' + '"'
$ws.Range("E13").Copy()
$ws.Range("E13").PasteSpecial(-4163)

$ws.Range("D14").Formula = '="' + '3.851.07' + '"'
$ws.Range("D14").Copy()
$ws.Range("D14").PasteSpecial(-4163)

$ws.Range("E14").Formula = '="' + '  +0.75%  ' + '"'
$ws.Range("E14").Copy()
$ws.Range("E14").PasteSpecial(-4163)

$ws.Range("E15").Formula = '="' + '  -2.14%  ' + '"'
$ws.Range("E15").Copy()
$ws.Range("E15").PasteSpecial(-4163)

$ws.Range("D16").Formula = '="' + '595.00' + '"'
$ws.Range("D16").Copy()
$ws.Range("D16").PasteSpecial(-4163)

$ws.Range("E16").Formula = '="' + '  -8.22%  ' + '"'
$ws.Range("E16").Copy()
$ws.Range("E16").PasteSpecial(-4163)

$ws.Range("D17").Formula = '="' + '66.189.04' + '"'
$ws.Range("D17").Copy()
$ws.Range("D17").PasteSpecial(-4163)

$ws.Range("E17").Formula = '="' + '  +0.66%  ' + '"'
$ws.Range("E17").Copy()
$ws.Range("E17").PasteSpecial(-4163)

$ws.Range("E18").Formula = '="' + '  +0.40%  ' + '"'
$ws.Range("E18").Copy()
$ws.Range("E18").PasteSpecial(-4163)

$ws.Range("D19").Formula = '="' + '3.320.86' + '"'
$ws.Range("D19").Copy()
$ws.Range("D19").PasteSpecial(-4163)

$ws.Range("E19").Formula = '="' + '  +0.82%  ' + '"'
$ws.Range("E19").Copy()
$ws.Range("E19").PasteSpecial(-4163)

$ws.Range("D20").Formula = '="' + '17.72' + '"'
$ws.Range("D20").Copy()
$ws.Range("D20").PasteSpecial(-4163)

$ws.Range("E20").Formula = '="' + '  -2.22%  ' + '"'
$ws.Range("E20").Copy()
$ws.Range("E20").PasteSpecial(-4163)

$ws.Range("D21").Formula = '="' + '10.91' + '"'
$ws.Range("D21").Copy()
$ws.Range("D21").PasteSpecial(-4163)

$ws.Range("E21").Formula = '="' + '  -4.31%  ' + '"'
$ws.Range("E21").Copy()
$ws.Range("E21").PasteSpecial(-4163)

$ws.Range("D22").Formula = '="' + '0.897' + '"'
$ws.Range("D22").Copy()
$ws.Range("D22").PasteSpecial(-4163)

$ws.Range("E22").Formula = '="' + '  -0.87%  ' + '"'
$ws.Range("E22").Copy()
$ws.Range("E22").PasteSpecial(-4163)

$ws.Range("D23").Formula = '="' + '17.92' + '"'
$ws.Range("D23").Copy()
$ws.Range("D23").PasteSpecial(-4163)

$ws.Range("E23").Formula = '="' + '  -1.13%  ' + '"'
$ws.Range("E23").Copy()
$ws.Range("E23").PasteSpecial(-4163)

$ws.Range("D24").Formula = '="' + '4.99' + '"'
$ws.Range("D24").Copy()
$ws.Range("D24").PasteSpecial(-4163)

$ws.Range("E24").Formula = '="' + '  +1.96%  ' + '"'
$ws.Range("E24").Copy()
$ws.Range("E24").PasteSpecial(-4163)

$ws.Range("D25").Formula = '="' + '98.86' + '"'
$ws.Range("D25").Copy()
$ws.Range("D25").PasteSpecial(-4163)

$ws.Range("E25").Formula = '="' + '  -7.95%  ' + '"'
$ws.Range("E25").Copy()
$ws.Range("E25").PasteSpecial(-4163)

$ws.Range("D26").Formula = '="' + '3.98' + '"'
$ws.Range("D26").Copy()
$ws.Range("D26").PasteSpecial(-4163)

$ws.Range("E26").Formula = '="' + '  +0.22%  ' + '"'
$ws.Range("E26").Copy()
$ws.Range("E26").PasteSpecial(-4163)

$ws.Range("E27").Formula = '="' + '  +1.14%  ' + '"'
$ws.Range("E27").Copy()
$ws.Range("E27").PasteSpecial(-4163)

$ws.Range("E28").Formula = '="' + '  -1.80%  ' + '"'
$ws.Range("E28").Copy()
$ws.Range("E28").PasteSpecial(-4163)

$ws.Range("E29").Formula = '="' + '  -2.46%  ' + '"'
$ws.Range("E29").Copy()
$ws.Range("E29").PasteSpecial(-4163)

$ws.Range("D30").Formula = '="' + '30.72' + '"'
$ws.Range("D30").Copy()
$ws.Range("D30").PasteSpecial(-4163)

$ws.Range("E30").Formula = '="' + '  +1.51%  ' + '"'
$ws.Range("E30").Copy()
$ws.Range("E30").PasteSpecial(-4163)

$ws.Range("D31").Formula = '="' + '6.65' + '"'
$ws.Range("D31").Copy()
$ws.Range("D31").PasteSpecial(-4163)

$ws.Range("E31").Formula = '="' + '  +5.88%  ' + '"'
$ws.Range("E31").Copy()
$ws.Range("E31").PasteSpecial(-4163)

$ws.Range("B32").Value = 'dogwifhat'

$ws.Range("C32").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'

$ws.Range("D32").Formula = '="' + '3.66' + '"'
$ws.Range("D32").Copy()
$ws.Range("D32").PasteSpecial(-4163)

$ws.Range("E32").Formula = '="' + '  -5.91%  ' + '"'
$ws.Range("E32").Copy()
$ws.Range("E32").PasteSpecial(-4163)

$ws.Range("B33").Value = 'Bittensor'

$ws.Range("C33").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'

$ws.Range("D33").Formula = '="' + '559.35' + '"'
$ws.Range("D33").Copy()
$ws.Range("D33").PasteSpecial(-4163)

$ws.Range("E33").Formula = '="' + '  +6.84%  ' + '"'
$ws.Range("E33").Copy()
$ws.Range("E33").PasteSpecial(-4163)

$ws.Range("D34").Formula = '="' + '10.87' + '"'
$ws.Range("D34").Copy()
$ws.Range("D34").PasteSpecial(-4163)

$ws.Range("E34").Formula = '="' + '  -1.50%  ' + '"'
$ws.Range("E34").Copy()
$ws.Range("E34").PasteSpecial(-4163)

$ws.Range("D35").Formula = '="' + '3.813.95' + '"'
$ws.Range("D35").Copy()
$ws.Range("D35").PasteSpecial(-4163)

$ws.Range("E35").Formula = '="' + '  +0.63%  ' + '"'
$ws.Range("E35").Copy()
$ws.Range("E35").PasteSpecial(-4163)

$ws.Range("E36").Formula = '="' + '  -1.35%  ' + '"'
$ws.Range("E36").Copy()
$ws.Range("E36").PasteSpecial(-4163)

$ws.Range("E37").Formula = '="' + '  +0.03%  ' + '"'
$ws.Range("E37").Copy()
$ws.Range("E37").PasteSpecial(-4163)

$ws.Range("E38").Formula = '="' + '  -2.68%  ' + '"'
$ws.Range("E38").Copy()
$ws.Range("E38").PasteSpecial(-4163)

$ws.Range("D39").Formula = '="' + '33.27' + '"'
$ws.Range("D39").Copy()
$ws.Range("D39").PasteSpecial(-4163)

$ws.Range("E39").Formula = '="' + '  +1.01%  ' + '"'
$ws.Range("E39").Copy()
$ws.Range("E39").PasteSpecial(-4163)

$ws.Range("D40").Formula = '="' + '0.128' + '"'
$ws.Range("D40").Copy()
$ws.Range("D40").PasteSpecial(-4163)

$ws.Range("E40").Formula = '="' + '  -1.36%  ' + '"'
$ws.Range("E40").Copy()
$ws.Range("E40").PasteSpecial(-4163)

$ws.Range("D41").Formula = '="' + '0.0' + [char]0x2083 + '0687' + '"'
$ws.Range("D41").Copy()
$ws.Range("D41").PasteSpecial(-4163)

$ws.Range("E41").Formula = '="' + '  -6.39%  ' + '"'
$ws.Range("E41").Copy()
$ws.Range("E41").PasteSpecial(-4163)

$ws.Range("E42").Formula = '="' + '  -6.08%  ' + '"'
$ws.Range("E42").Copy()
$ws.Range("E42").PasteSpecial(-4163)

$ws.Range("D43").Formula = '="' + '3.38' + '"'
$ws.Range("D43").Copy()
$ws.Range("D43").PasteSpecial(-4163)

$ws.Range("E43").Formula = '="' + '  +3.95%  ' + '"'
$ws.Range("E43").Copy()
$ws.Range("E43").PasteSpecial(-4163)

$ws.Range("E44").Formula = '="' + '  -4.14%  ' + '"'
$ws.Range("E44").Copy()
$ws.Range("E44").PasteSpecial(-4163)

$ws.Range("D45").Formula = '="' + '0.332' + '"'
$ws.Range("D45").Copy()
$ws.Range("D45").PasteSpecial(-4163)

$ws.Range("E46").Formula = '="' + '  -0.85%  ' + '"'
$ws.Range("E46").Copy()
$ws.Range("E46").PasteSpecial(-4163)

$ws.Range("D47").Formula = '="' + '3.06' + '"'
$ws.Range("D47").Copy()
$ws.Range("D47").PasteSpecial(-4163)

$ws.Range("E47").Formula = '="' + '  -8.36%  ' + '"'
$ws.Range("E47").Copy()
$ws.Range("E47").PasteSpecial(-4163)

$ws.Range("E48").Formula = '="' + '  -2.37%  ' + '"'
$ws.Range("E48").Copy()
$ws.Range("E48").PasteSpecial(-4163)

$ws.Range("E49").Formula = '="' + '  +0.15%  ' + '"'
$ws.Range("E49").Copy()
$ws.Range("E49").PasteSpecial(-4163)

$ws.Range("E50").Formula = '="' + '  -2.50%  ' + '"'
$ws.Range("E50").Copy()
$ws.Range("E50").PasteSpecial(-4163)

$ws.Range("D51").Formula = '="' + '129.57' + '"'
$ws.Range("D51").Copy()
$ws.Range("D51").PasteSpecial(-4163)

$ws.Range("E51").Formula = '="' + '  +6.21%  ' + '"'
$ws.Range("E51").Copy()
$ws.Range("E51").PasteSpecial(-4163)

$excel.CutCopyMode = 0
